{"js": "// Fixed player running into VGA HEX bug. Added the corresponding log entry.\n// Insert two new log paragraphs (style \"Normal1\") right after the paragraph\n// that ends the existing \"VGA HEX decoder module\" log entry.\n\nconst body = context.document.body;\n\n// Locate the paragraph to anchor on via its unique, stable text.\nconst results = body.search(\"VGA HEX decoder module\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Anchor paragraph ('VGA HEX decoder module') not found.\");\n}\n\nconst anchorParagraph = results.items[0].paragraphs.getFirst();\n\n// First new paragraph, inserted directly after the anchor paragraph.\nconst newParagraph1 = anchorParagraph.insertParagraph(\n  \"             - Fixed a bug where player moving left will loop into the VGA HEX panels if player\\u2019\",\n  \"After\"\n);\nnewParagraph1.style = \"Normal1\";\n\n// Second new paragraph, inserted directly after the first new paragraph.\nconst newParagraph2 = newParagraph1.insertParagraph(\n  \"y coordinate is within the y range of the VGA HEX panels.\",\n  \"After\"\n);\nnewParagraph2.style = \"Normal1\";\n\nawait context.sync();\n", "ps1": "# Fixed player running into VGA HEX bug. Added the corresponding log entry.\n# Insert two new log paragraphs (style \"Normal1\") right after the paragraph\n# that ends the existing \"VGA HEX decoder module\" log entry.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph via its unique, stable text.\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"VGA HEX decoder module\")\nif (-not $found) {\n    throw \"Anchor paragraph ('VGA HEX decoder module') not found.\"\n}\n$anchorParagraph = $searchRange.Paragraphs(1)\n\n# First new paragraph, inserted directly after the anchor paragraph.\n$anchorParagraph.Range.InsertParagraphAfter()\n$newParagraph1 = $anchorParagraph.Next()\n$newParagraph1.Range.Text = \"             - Fixed a bug where player moving left will loop into the VGA HEX panels if player\" + [char]0x2019\n$newParagraph1.Style = \"Normal1\"\n\n# Second new paragraph, inserted directly after the first new paragraph.\n$newParagraph1.Range.InsertParagraphAfter()\n$newParagraph2 = $newParagraph1.Next()\n$newParagraph2.Range.Text = \"y coordinate is within the y range of the VGA HEX panels.\"\n$newParagraph2.Style = \"Normal1\"\n"}
